# "quadro de lideranças adicionado"
# Adds two new worksheets (Planilha2, Planilha3) after the existing
# Planilha1, fills them with the leadership-board data, and leaves
# Planilha3 as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the two new sheets, in order, right after Planilha1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Planilha2"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Planilha3"

# --- Planilha2: Diretores ---
$ws2.Range("A2").Value = "Diretores"
$ws2.Range("B1").Value = "Presidente"
$ws2.Range("C1").Value = "Marketing"
$ws2.Range("D1").Value = "Comercial"
$ws2.Range("E1").Value = "Projetos"
$ws2.Range("F1").Value = "Gente_e_Gestao"
$ws2.Range("B2").Value = "Daniel Fonseca"
$ws2.Range("C2").Value = "Aline Maia"
$ws2.Range("D2").Value = "Vinius Santi"
$ws2.Range("E2").Value = "Jardel Salles"
$ws2.Range("F2").Value = "Kleber Azeredo"

# --- Planilha3: Coordenadores ---
$ws3.Range("A2").Value = "Coordenadores"
$ws3.Range("C1").Value = "Comp"
$ws3.Range("D1").Value = "Eletr"
$ws3.Range("E1").Value = "Autom"
$ws3.Range("F1").Value = "Endo"
$ws3.Range("G1").Value = "Inbound"
$ws3.Range("H1").Value = "PosV"
$ws3.Range("I1").Value = "Vendas"
$ws3.Range("J1").Value = "Gente"
$ws3.Range("K1").Value = "Gestao"
$ws3.Range("B2").Value = "Rayssa Alves"
$ws3.Range("C2").Value = "Higor Brandão"
$ws3.Range("D2").Value = "Laura Barros"
$ws3.Range("E2").Value = "Diogo Ribeiro"
$ws3.Range("F2").Value = "Guilherme Barreto"
$ws3.Range("G2").Value = "Micaella Barcellos"
$ws3.Range("H2").Value = "João Carvalho"
$ws3.Range("I2").Value = "Matheus Henrique"
$ws3.Range("J2").Value = "Talita Silva"
$ws3.Range("K2").Value = "Giana Bastos"
$ws3.Range("B1").Value = "ArqEUrb"

# --- restore per-sheet selections; select in tab order so Planilha3
#     (the last one touched) ends up the active/visible tab ---
$ws1.Range("F7").Select() | Out-Null
$ws2.Range("C5").Select() | Out-Null
$ws3.Range("B1").Select() | Out-Null
